$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the histogram delimiter values (channel order shifted)
$ws.Range("A3").Value = 465
$ws.Range("B3").Value = 805
$ws.Range("C3").Value = 905
$ws.Range("D3").Value = 1160

$ws.Range("A4").Value = 818
$ws.Range("B4").Value = 1153

$ws.Range("A5").Value = 1176
$ws.Range("B5").Value = 1505

$ws.Range("A6").Value = 1516
$ws.Range("B6").Value = 1871

$ws.Range("A7").Value = 1888
$ws.Range("B7").Value = 2233

# Update the active selection to reflect where the user last clicked
$ws.Range("C10").Select()
